$d = $word.ActiveDocument

# Remove every paragraph after the first one (the "Amy Han Hsun Shih" title
# in Heading1 style). This deletes the subtitle, date, salutation and the
# body paragraph(s) that follow, leaving only the title paragraph before
# the section properties.
$paraCount = $d.Paragraphs.Count
if ($paraCount -gt 1) {
    $startPara = $d.Paragraphs.Item(2)
    $endPara = $d.Paragraphs.Item($paraCount)
    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}
